$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.615.33'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '3.466.03'
$ws.Range('E3').Value = '  +2.95%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''573.44'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').Value = '''159.73'
$ws.Range('E6').Value = '  +4.30%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.477.30'
$ws.Range('E8').Value = '  +3.08%  '
$ws.Range('D9').Value = '''0.583'
$ws.Range('E9').Value = '  +11.26%  '
$ws.Range('D10').Value = '''7.36'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').Value = '''0.125'
$ws.Range('E11').Value = '  +5.33%  '
$ws.Range('D12').Value = '''0.445'
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').Value = '4.058.94'
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '''0.0000194'
$ws.Range('E15').Value = '  +7.68%  '
$ws.Range('D16').Value = '''28.32'
$ws.Range('E16').Value = '  +5.10%  '
$ws.Range('D17').Value = '64.675.75'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '3.488.90'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').Value = '''6.42'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = '''14.36'
$ws.Range('E20').Value = '  +3.39%  '
$ws.Range('D21').Value = '''391.38'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = '''8.25'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('D23').Value = '''73.56'
$ws.Range('E23').Value = '  +4.60%  '
$ws.Range('D24').Value = '''0.543'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '''0.0000120'
$ws.Range('E26').Value = '  +24.17%  '
$ws.Range('D27').Value = '''9.55'
$ws.Range('E27').Value = '  +2.04%  '
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '''1.01'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '''1.46'
$ws.Range('E30').Value = '  +11.88%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''6.15'
$ws.Range('E31').Value = '  +10.34%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').Value = '''6.56'
$ws.Range('E33').Value = '  +3.60%  '
$ws.Range('D34').Value = '''23.64'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '''7.00'
$ws.Range('E36').Value = '  +4.19%  '
$ws.Range('D37').Value = '''1.49'
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').Value = '''160.86'
$ws.Range('E38').Value = '  +2.03%  '
$ws.Range('D39').Value = '''1.89'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').Value = '''0.0774'
$ws.Range('E40').Value = '  +4.30%  '
$ws.Range('D41').Value = '''27.48'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('D42').Value = '2.924.53'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('D43').Value = '''0.0318'
$ws.Range('E43').Value = '  -2.69%  '
$ws.Range('D44').Value = '''42.65'
$ws.Range('E44').Value = '  +4.42%  '
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('D46').Value = '''0.774'
$ws.Range('E46').Value = '  +3.38%  '
$ws.Range('D47').Value = '''23.74'
$ws.Range('E47').Value = '  +8.27%  '
$ws.Range('E48').Value = '  +6.43%  '
$ws.Range('D49').Value = '''2.22'
$ws.Range('E49').Value = '  +21.05%  '
$ws.Range('E50').Value = '  +4.79%  '
$ws.Range('D51').Value = '''0.856'
$ws.Range('E51').Value = '  +6.54%  '
